# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (fund-holdings detail) right after "2021-Q4".
# 2) Rebuild the "总计" (totals) summary sheet with a new leading row for
#    2022-Q1, shifting the previous rows down by one.
#
# The totals sheet is dropped and re-added so the new sheetId ordering
# matches the source data (2022-Q1 gets sheetId 5, 总计 becomes sheetId 6).

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotals = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$oldTotals.Delete()

# ---------------------------------------------------------------------------
# New detail sheet: "2022-Q1"
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Reuse the header / index-column look (bold, centered, bordered) from the
# existing "2021-Q4" sheet instead of re-creating a style from scratch.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A5").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3

# Columns B (fund code) and D:G (numeric-looking figures kept as text in the
# source data) must stay text, otherwise leading zeros / trailing zeros are
# lost to Excel's automatic number coercion.
$q1.Range("B2:B5").NumberFormat = "@"
$q1.Range("D2:G5").NumberFormat = "@"

$q1.Range("B2").Value = "070021"
$q1.Range("C2").Value = "嘉实主题新动力混合"
$q1.Range("D2").Value = "24.04"
$q1.Range("E2").Value = "93.93"
$q1.Range("F2").Value = "4.57"
$q1.Range("G2").Value = "1.0986"
$q1.Range("H2").Value = 9

$q1.Range("B3").Value = "000985"
$q1.Range("C3").Value = "嘉实逆向策略股票"
$q1.Range("D3").Value = "13.64"
$q1.Range("E3").Value = "93.90"
$q1.Range("F3").Value = "4.58"
$q1.Range("G3").Value = "0.6247"
$q1.Range("H3").Value = 9

$q1.Range("B4").Value = "014307"
$q1.Range("C4").Value = "嘉实多元动力混合A"
$q1.Range("D4").Value = "1.83"
$q1.Range("E4").Value = "91.81"
$q1.Range("F4").Value = "4.27"
$q1.Range("G4").Value = "0.0781"
$q1.Range("H4").Value = 9

$q1.Range("B5").Value = "014308"
$q1.Range("C5").Value = "嘉实多元动力混合C"
$q1.Range("D5").Value = "0.20"
$q1.Range("E5").Value = "91.81"
$q1.Range("F5").Value = "4.27"
$q1.Range("G5").Value = "0.0085"
$q1.Range("H5").Value = 9

# ---------------------------------------------------------------------------
# Rebuilt totals sheet: "总计"
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Add($null, $q1)
$totals.Name = "总计"

$q1.Range("B1:D1").Copy()
$totals.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$totals.Range("A2:A6").PasteSpecial(-4122)

$totals.Range("B1").Value = "日期"
$totals.Range("C1").Value = "持有数量(只)"
$totals.Range("D1").Value = "持有市值(亿元)"

$totals.Range("A2").Value = 0
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4

$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 4
$totals.Range("D2").Value = 1.81

$totals.Range("B3").Value = "2021-Q4"
$totals.Range("C3").Value = 8
$totals.Range("D3").Value = 2.87

$totals.Range("B4").Value = "2021-Q3"
$totals.Range("C4").Value = 4
$totals.Range("D4").Value = 1.51

$totals.Range("B5").Value = "2021-Q2"
$totals.Range("C5").Value = 5
$totals.Range("D5").Value = 0.33

$totals.Range("B6").Value = "2021-Q1"
$totals.Range("C6").Value = 3
$totals.Range("D6").Value = 0.42
